$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Email " (with trailing space) -> "Email" (no trailing space)
$ws.Range("C1").Value = "Email"

# Row 4 (Sally) gets a blank Email cell, styled like the hyperlink cells above it
# (no hyperlink actually added, matching the source data)
$ws.Range("C4").Style = "Hyperlink"

# Minor column width tweaks
$ws.Columns.Item(3).ColumnWidth = 19
$ws.Columns.Item(5).ColumnWidth = 9.67

# Update the active selection to C4
$ws.Range("C4").Select()
